$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2-26 of Sheet1, columns A-T (20 cols), reflecting the refreshed TPM-derived NATMI stats
$data = @(
    @("ECs", "Ntrk3", "Ptprs", "ECs", 1, 0.3333333333333333, 0.002145666666666667, 0.006437, 0.00807320947388686, 0.008843624333499573, 3, 1, 2.005664333333333, 6.016992999999999, 0.02976781902817159, 0.03172257287647481, 0.004303487104555556, 0.03873138394099999, 0.0002403218385951844, 0.0002805425174116062),
    @("ECs", "Ntrk3", "Ptprs", "FAPs", 1, 0.3333333333333333, 0.002145666666666667, 0.006437, 0.00807320947388686, 0.008843624333499573, 3, 1, 32.65736066666667, 97.972082, 0.4846964599741412, 0.5165248673390457, 0.07007181020377778, 0.630646291834, 0.00391305605262266, 0.004567951885657224),
    @("ECs", "Ntrk3", "Ptprs", "Inflammatory-Mac", 1, 0.3333333333333333, 0.002145666666666667, 0.006437, 0.00807320947388686, 0.008843624333499573, 3, 1, 10.495667, 31.487001, 0.1557753760903259, 0.1660046278737794, 0.02252020282633334, 0.202681825437, 0.001257607242050707, 0.001468082566538097),
    @("ECs", "Ntrk3", "Ptprs", "MuSCs", 1, 0.3333333333333333, 0.002145666666666667, 0.006437, 0.00807320947388686, 0.008843624333499573, 2, 1, 12.4553565, 24.910713, 0.1848608423958749, 0.1313333601264699, 0.0267250432635, 0.160350259581, 0.001492420304181083, 0.001161462899414712),
    @("ECs", "Ntrk3", "Ptprs", "Resolving-Mac", 1, 0.3333333333333333, 0.002145666666666667, 0.006437, 0.00807320947388686, 0.008843624333499573, 3, 1, 9.762884, 29.288652, 0.1448995025114864, 0.1544145717842301, 0.02094789476933333, 0.188531052924, 0.001169804036437225, 0.001365584464477934),
    @("FAPs", "Ntrk3", "Ptprs", "ECs", 3, 1, 0.03558433333333334, 0.106753, 0.1338883534202026, 0.1466651279282398, 3, 1, 2.005664333333333, 6.016992999999999, 0.02976781902817159, 0.03172257287647481, 0.07137022819211111, 0.6423320537289999, 0.003985564274592469, 0.00465259520914109),
    @("FAPs", "Ntrk3", "Ptprs", "FAPs", 3, 1, 0.03558433333333334, 0.106753, 0.1338883534202026, 0.1466651279282398, 3, 1, 32.65736066666667, 97.972082, 0.4846964599741412, 0.5165248673390457, 1.162090407749556, 10.458813669746, 0.0648952109345389, 0.07575618574639824),
    @("FAPs", "Ntrk3", "Ptprs", "Inflammatory-Mac", 3, 1, 0.03558433333333334, 0.106753, 0.1338883534202026, 0.1466651279282398, 3, 1, 10.495667, 31.487001, 0.1557753760903259, 0.1660046278737794, 0.3734813130836667, 3.361331817753, 0.02085650860814652, 0.0243470899837877),
    @("FAPs", "Ntrk3", "Ptprs", "MuSCs", 3, 1, 0.03558433333333334, 0.106753, 0.1338883534202026, 0.1466651279282398, 2, 1, 12.4553565, 24.910713, 0.1848608423958749, 0.1313333601264699, 0.4432155574815, 2.659293344889, 0.02475071380025526, 0.0192620240641943),
    @("FAPs", "Ntrk3", "Ptprs", "Resolving-Mac", 3, 1, 0.03558433333333334, 0.106753, 0.1338883534202026, 0.1466651279282398, 3, 1, 9.762884, 29.288652, 0.1448995025114864, 0.1544145717842301, 0.3474057185506667, 3.126651466956, 0.01940035580266942, 0.02264723292471848),
    @("Inflammatory-Mac", "Ntrk3", "Ptprs", "ECs", 1, 0.3333333333333333, 0.04573, 0.13719, 0.172062079807758, 0.188481718551003, 3, 1, 2.005664333333333, 6.016992999999999, 0.02976781902817159, 0.03172257287647481, 0.09171902996333332, 0.82547126967, 0.005121912853328157, 0.005979125052617408),
    @("Inflammatory-Mac", "Ntrk3", "Ptprs", "FAPs", 1, 0.3333333333333333, 0.04573, 0.13719, 0.172062079807758, 0.188481718551003, 3, 1, 32.65736066666667, 97.972082, 0.4846964599741412, 0.5165248673390457, 1.493421103286667, 13.44078992958, 0.08339788097860848, 0.09735549467039217),
    @("Inflammatory-Mac", "Ntrk3", "Ptprs", "Inflammatory-Mac", 1, 0.3333333333333333, 0.04573, 0.13719, 0.172062079807758, 0.188481718551003, 3, 1, 10.495667, 31.487001, 0.1557753760903259, 0.1660046278737794, 0.47996685191, 4.31970166719, 0.02680303519293717, 0.03128883754906968),
    @("Inflammatory-Mac", "Ntrk3", "Ptprs", "MuSCs", 1, 0.3333333333333333, 0.04573, 0.13719, 0.172062079807758, 0.188481718551003, 2, 1, 12.4553565, 24.910713, 0.1848608423958749, 0.1313333601264699, 0.569583452745, 3.41750071647, 0.03180754101764839, 0.02475393741971482),
    @("Inflammatory-Mac", "Ntrk3", "Ptprs", "Resolving-Mac", 1, 0.3333333333333333, 0.04573, 0.13719, 0.172062079807758, 0.188481718551003, 3, 1, 9.762884, 29.288652, 0.1448995025114864, 0.1544145717842301, 0.44645668532, 4.01811016788, 0.0249317097652358, 0.0291043238592089),
    @("MuSCs", "Ntrk3", "Ptprs", "ECs", 2, 1, 0.06945950000000001, 0.138919, 0.2613458568206203, 0.1908571459974254, 3, 1, 2.005664333333333, 6.016992999999999, 0.02976781902817159, 0.03172257287647481, 0.1393124417611667, 0.835874650567, 0.00777969616959867, 0.00605447972289932),
    @("MuSCs", "Ntrk3", "Ptprs", "FAPs", 2, 1, 0.06945950000000001, 0.138919, 0.2613458568206203, 0.1908571459974254, 3, 1, 32.65736066666667, 97.972082, 0.4846964599741412, 0.5165248673390457, 2.268363943226334, 13.610183659358, 0.1266734116298635, 0.09858246201702901),
    @("MuSCs", "Ntrk3", "Ptprs", "Inflammatory-Mac", 2, 1, 0.06945950000000001, 0.138919, 0.2613458568206203, 0.1908571459974254, 3, 1, 10.495667, 31.487001, 0.1557753760903259, 0.1660046278737794, 0.7290237819865002, 4.374142691919001, 0.04071124913588058, 0.03168316949835418),
    @("MuSCs", "Ntrk3", "Ptprs", "MuSCs", 2, 1, 0.06945950000000001, 0.138919, 0.2613458568206203, 0.1908571459974254, 2, 1, 12.4553565, 24.910713, 0.1848608423958749, 0.1313333601264699, 0.8651428348117501, 3.460571339247001, 0.04831261524853157, 0.02506591028799011),
    @("MuSCs", "Ntrk3", "Ptprs", "Resolving-Mac", 2, 1, 0.06945950000000001, 0.138919, 0.2613458568206203, 0.1908571459974254, 3, 1, 9.762884, 29.288652, 0.1448995025114864, 0.1544145717842301, 0.6781250411980001, 4.068750247188, 0.03786888463674604, 0.02947112447115272),
    @("Resolving-Mac", "Ntrk3", "Ptprs", "ECs", 1, 0.3333333333333333, 0.1128566666666667, 0.33857, 0.4246305004775321, 0.4651523831898322, 3, 1, 2.005664333333333, 6.016992999999999, 0.02976781902817159, 0.03172257287647481, 0.2263525911122222, 2.03717332001, 0.01264032389205711, 0.01475583037440539),
    @("Resolving-Mac", "Ntrk3", "Ptprs", "FAPs", 1, 0.3333333333333333, 0.1128566666666667, 0.33857, 0.4246305004775321, 0.4651523831898322, 3, 1, 32.65736066666667, 97.972082, 0.4846964599741412, 0.5165248673390457, 3.685600866971111, 33.17040780274, 0.2058169003785077, 0.240262773019569),
    @("Resolving-Mac", "Ntrk3", "Ptprs", "Inflammatory-Mac", 1, 0.3333333333333333, 0.1128566666666667, 0.33857, 0.4246305004775321, 0.4651523831898322, 3, 1, 10.495667, 31.487001, 0.1557753760903259, 0.1660046278737794, 1.184505992063333, 10.66055392857, 0.06614697591131086, 0.07721744827602973),
    @("Resolving-Mac", "Ntrk3", "Ptprs", "MuSCs", 1, 0.3333333333333333, 0.1128566666666667, 0.33857, 0.4246305004775321, 0.4651523831898322, 2, 1, 12.4553565, 24.910713, 0.1848608423958749, 0.1313333601264699, 1.405670016735, 8.43402010041, 0.07849755202525853, 0.06109002545515595),
    @("Resolving-Mac", "Ntrk3", "Ptprs", "Resolving-Mac", 1, 0.3333333333333333, 0.1128566666666667, 0.33857, 0.4246305004775321, 0.4651523831898322, 3, 1, 9.762884, 29.288652, 0.1448995025114864, 0.1544145717842301, 1.101806545293333, 9.91625890764, 0.06152874827039787, 0.07182630606467204)
)

$startRow = 2
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $data[$i]
    $r = $startRow + $i
    for ($j = 0; $j -lt $row.Length; $j++) {
        $ws.Cells.Item($r, $j + 1).Value = $row[$j]
    }
}
